# Repull data, push all data, mean calculation
# Update the dSF (F column) values to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 3
    10 = -1
    16 = 0
    18 = 2
    20 = 4
    31 = 3
    41 = -1
    43 = 1
    45 = -1
    46 = -1
    47 = -5
    55 = 1
    57 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
